# Agencies.xlsx edit script
# Renames the sole "Canada Space Agency" sheet to "CSA" and adds two new
# sheets ("Roscosmos" and "ISRO") populated with agency data, then makes
# "Roscosmos" the active tab.

$wb = $excel.ActiveWorkbook
$csa = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Create the three sheets we need while carefully controlling the
# internal sheetId counter (it is "max used id + 1" and is not reused
# when a sheet is removed from the middle, only when the highest id is
# removed). The target workbook requires:
#   CSA        sheetId=1 (unchanged, original sheet)
#   Roscosmos  sheetId=4
#   ISRO       sheetId=2
# so we create ISRO first (id 2), then two scratch sheets (ids 3 and 4),
# delete the id-4 scratch sheet, create Roscosmos (which then reclaims
# id 4), and finally delete the leftover id-3 scratch sheet.
# ---------------------------------------------------------------------

# 1) Copy CSA right after itself -> position 2, sheetId=2 (future ISRO)
$csa.Copy($null, $csa)

# 2) Copy position 2 right after itself -> position 3, sheetId=3 (scratch)
$wb.Worksheets.Item(2).Copy($null, $wb.Worksheets.Item(2))

# 3) Copy position 3 right after itself -> position 4, sheetId=4 (scratch)
$wb.Worksheets.Item(3).Copy($null, $wb.Worksheets.Item(3))

# Fill in the ISRO data now (position 2) while the shared-string table
# only contains the original strings, so the new strings land in the
# same order as in the target file.
$isro = $wb.Worksheets.Item(2)
$isro.Range("A2").Value = 3
$isro.Range("B2").Value = 38571
$isro.Range("C2").Value = "India"
$isro.Range("D2").Value = 2
$isro.Range("E2").Value = 1360
$isro.Range("F2").Value = "Shira"
$isro.Range("G2").Value = 38968
$isro.Range("H2").Value = 41254
$isro.Range("I2").Value = "Completed"
$isro.Range("J2").Value = "Solar activity observation"
$isro.Range("F3").Value = "Sari"
$isro.Range("G3").Value = 42078
$isro.Range("H3").Clear()
$isro.Range("I3").Value = "Activa"
$isro.Range("J3").Value = "Observation of far Universe objects"

# 4) Remove the id-4 scratch sheet (now at position 4)
$wb.Worksheets.Item(4).Delete()

# 5) Copy CSA (position 1) right after itself -> position 2, reclaiming
#    sheetId=4 for the future Roscosmos sheet, pushing ISRO to position 3.
$wb.Worksheets.Item(1).Copy($null, $wb.Worksheets.Item(1))

$roscosmos = $wb.Worksheets.Item(2)
$roscosmos.Range("A2").Value = 4
$roscosmos.Range("B2").Value = 34611
$roscosmos.Range("C2").Value = "Russia"
$roscosmos.Range("D2").Value = 1.58
$roscosmos.Range("E2").Value = 144.5
$roscosmos.Range("F2").Value = "Proton"
$roscosmos.Range("G2").Value = 38538
$roscosmos.Range("H2").Clear()
$roscosmos.Range("I2").Value = "Active"
$roscosmos.Range("J2").Value = "ISS supplement and crew delivery"
$roscosmos.Range("F3").Value = "Vostok"
$roscosmos.Range("G3").Value = 37655
$roscosmos.Range("H3").Value = 39756
$roscosmos.Range("I3").Value = "Completed"
$roscosmos.Range("J3").Value = "Earth atmosphere observation"

# 6) Remove the id-3 scratch sheet (now pushed to position 4)
$wb.Worksheets.Item(4).Delete()

# ---------------------------------------------------------------------
# Rename sheets to their final names (position based, safe now that no
# more insertions/deletions happen).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "CSA"
$wb.Worksheets.Item(2).Name = "Roscosmos"
$wb.Worksheets.Item(3).Name = "ISRO"

# ---------------------------------------------------------------------
# Update the CSA sheet's own data: Population (E2) 12 -> 37.9, and the
# selection changes to a header-row range select.
# ---------------------------------------------------------------------
$csaFinal = $wb.Worksheets.Item(1)
$csaFinal.Range("E2").Value = 37.9
$csaFinal.Range("A1:J1").Select()

# Roscosmos keeps a single-cell selection on J3 and becomes the active tab.
$roscosmosFinal = $wb.Worksheets.Item(2)
$roscosmosFinal.Range("J3").Select()
$roscosmosFinal.Activate()

# ISRO keeps a single-cell selection on D2.
$isroFinal = $wb.Worksheets.Item(3)
$isroFinal.Range("D2").Select()

# Re-activate Roscosmos last so it ends up as the active/selected tab.
$roscosmosFinal.Activate()
